# Apply updated dSF (column F) values as per repull of data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 0
$ws.Range("F4").Value = -3
$ws.Range("F5").Value = -2
$ws.Range("F7").Value = -1
$ws.Range("F10").Value = -3
$ws.Range("F11").Value = -4
$ws.Range("F17").Value = -4
$ws.Range("F21").Value = -11
$ws.Range("F22").Value = -2
$ws.Range("F23").Value = -2
$ws.Range("F28").Value = -3
$ws.Range("F31").Value = -4
$ws.Range("F34").Value = 5
